$wb = $excel.ActiveWorkbook

# The edit targets the "Card13" sheet (the dimension A1:M13 -> A1:N13 sheet).
$ws = $wb.Worksheets.Item("Card13")

# --- 1. Add the new header cell N1 = "Event " -------------------------------
# Copy M1's formatting (bold header style with border/alignment) into N1
# before we touch M1's own text, then set N1's text.
$ws.Range("M1").Copy($ws.Range("N1"))
$ws.Range("N1").Value = "Event "

# --- 2. Create the new "Event" data cells N2:N13 -----------------------------
# These are blank placeholder cells (mirroring the blank "nan"-column cells
# elsewhere in the sheet before they are populated). A formula that evaluates
# to an empty string is used so the cell is actually persisted as a blank
# text cell (a literal empty value is dropped/un-persisted by the engine).
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 14).Formula = '=""'
}

# --- 3. Fix the M1 header text: "Correction " -> "Correction" ---------------
$ws.Range("M1").Value = "Correction"

# --- 4. Populate M2:M13 with "nan" (placeholder text, matching every other
#        blank cell in those rows) now that M1's header describes real data --
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 13).Value = "nan"
}
